# trabalho 07 - Bellman-Ford
# Updates the adjacency/weight matrix on sheet "grafo2" and leaves that
# sheet as the active one (matching the view state captured in the diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("grafo2")

# --- Data edits on grafo2 (weights changed for the Bellman-Ford exercise) ---
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = 0
$ws.Range("M2").Value = 0

$ws.Range("D3").Value = 0
$ws.Range("F3").Value = 1
$ws.Range("J3").Value = 0

$ws.Range("E4").Value = 0
$ws.Range("H4").Value = 0
$ws.Range("K4").Value = 0

$ws.Range("H5").Value = 4
$ws.Range("O5").Value = 0

$ws.Range("H7").Value = 1
$ws.Range("J7").Value = 0
$ws.Range("N7").Value = 0

$ws.Range("J9").Value = 10
$ws.Range("N9").Value = 0

$ws.Range("L11").Value = 0
$ws.Range("M11").Value = 0

$ws.Range("N12").Value = 9
$ws.Range("O12").Value = 6

# --- View state: make grafo2 the active sheet/tab ---
[void]$ws.Activate()

# Page setup matching the captured print settings for grafo2
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Scroll/zoom/selection state for grafo2's window
$excel.ActiveWindow.Zoom = 60
$excel.ActiveWindow.ScrollColumn = 6
$excel.ActiveWindow.ScrollRow = 1
[void]$ws.Range("O5").Select()

Write-Output "edit applied"
